# "Updated to do list"
# - Fill in row 34 (6th task under the "Server" section) with the new
#   task, owner, % complete and date.
# - Add a reviewer comment on B34.
# - Move the selection/viewport down to the newly edited row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task row (A34 already holds the task number "6").
$ws.Range("B34").Value = "Tidy up data formatting in Excel"
$ws.Range("C34").Value = "FC"
$ws.Range("D34").Value = 0.5
$ws.Range("E34").Value = 41094

# Leave a note for the task (legacy cell comment, same style Excel uses
# for Review > New Comment: bold author name followed by the note text).
$excel.UserName = "Fintan Costello"
$comment = $ws.Range("B34").AddComment("Fintan Costello:`ncontributions done, check with Eoin about votes")

# Move the active selection/view to where the edit happened.
$ws.Range("F34").Select()
